$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of test-case data to append below the existing header + row2.
# Columns: A = Test Case ID, B = Test scenario, C = Status
$rows = @(
    @(" iAU_TC_ID_239", "Exam_Prerequisit_for_iAU_TC_ID_239", "passed"),
    @("iAU_TC_ID_239", "Exam_Prerequisit_for_iAU_TC_ID_239", "passed"),
    @("iAU_TC_ID_239", "@RegressionA candidate answering questions", "passed"),
    @("iAU_TC_ID_239", "Exam_Prerequisit_for_iAU_TC_ID_239", "passed"),
    @("iAU_TC_ID_218", '@RegressionA "Validation of Delivery--> Marking"', "passed"),
    @("iAU_TC_ID_240.,iAU_TC_ID_241", "@RegressionA Verify Marker score card", "passed"),
    @("iAU_TC_ID_221.,iAU_TC_ID_222", "@RegressionA Verify Elumina", "passed"),
    @("iAU_TC_ID_242", "@RegressionA Verify Elumina Markers", "passed"),
    @("iAU_TC_ID_227", "@RegressionA Verify Elumina Markers", "failed")
)

# Row 2 already exists (iAU_TC_ID_108 test case); rewrite it and append rows 3-10.
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
